$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.15040000000001
$ws.Range("D6").Value = -7.836199999999998
$ws.Range("A14").Value = -21.74469999999998
$ws.Range("B15").Value = 5.125099999999999
$ws.Range("A16").Value = -22.10240000000001
$ws.Range("D18").Value = -8.587499999999999
$ws.Range("D19").Value = -8.235899999999997
$ws.Range("A21").Value = -21.77649999999998
$ws.Range("B21").Value = 5.919
$ws.Range("B22").Value = 9.502700000000001
$ws.Range("A23").Value = -20.33999999999999
$ws.Range("B24").Value = 5.550800000000002
$ws.Range("A25").Value = -21.72549999999998
$ws.Range("A26").Value = -21.08919999999996
$ws.Range("B27").Value = 5.136800000000002
$ws.Range("B28").Value = 4.906600000000001
$ws.Range("A29").Value = -21.67389999999997
$ws.Range("D35").Value = -8.259799999999995
$ws.Range("B36").Value = 9.3804
$ws.Range("B39").Value = 9.152800000000003
$ws.Range("A40").Value = -20.758
$ws.Range("D44").Value = -7.171999999999998
$ws.Range("B45").Value = 5.028700000000003
$ws.Range("D47").Value = -7.175000000000009
$ws.Range("B48").Value = 7.658100000000003
$ws.Range("B49").Value = 6.037099999999997
$ws.Range("D50").Value = -8.255300000000004
$ws.Range("D51").Value = -8.128299999999998
$ws.Range("B52").Value = 4.679200000000001
$ws.Range("D52").Value = -8.268599999999998
$ws.Range("A53").Value = -21.61209999999998
$ws.Range("B53").Value = 6.124999999999996
$ws.Range("B54").Value = 4.616700000000004
$ws.Range("D55").Value = -8.889600000000002
$ws.Range("A57").Value = -22.52800000000002
$ws.Range("B57").Value = 4.437199999999997
$ws.Range("D57").Value = -8.262499999999999
$ws.Range("D58").Value = -8.069500000000012
$ws.Range("A59").Value = -22.64630000000001
$ws.Range("D64").Value = -7.493499999999992
$ws.Range("A65").Value = -21.83669999999998
$ws.Range("D66").Value = -7.555599999999991
$ws.Range("A69").Value = -21.6245
$ws.Range("B70").Value = 8.371700000000004
$ws.Range("B71").Value = 5.025299999999996
$ws.Range("A79").Value = -20.50720000000001
$ws.Range("D80").Value = -7.613099999999999
$ws.Range("A83").Value = -21.59139999999999
$ws.Range("D83").Value = -7.249100000000006
$ws.Range("B86").Value = 4.801200000000006
$ws.Range("B87").Value = 5.222900000000002
$ws.Range("B89").Value = 4.4565
$ws.Range("A91").Value = -20.94829999999999
$ws.Range("D92").Value = -6.681500000000001
$ws.Range("A93").Value = -21.49690000000002
$ws.Range("D94").Value = -6.670000000000003
$ws.Range("D96").Value = -8.519400000000003
$ws.Range("D97").Value = -8.618799999999995
$ws.Range("A100").Value = -22.1057
$ws.Range("B101").Value = 4.941500000000001
$ws.Range("D101").Value = -7.723299999999998
$ws.Range("A103").Value = -21.72710000000001
